$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("SignIn"): duplicate row 2 into a new row 3 with updated
# email / password values, then rebuild the hyperlinks on C2:D3.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SignIn")

$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))

$ws1.Range("C2").Value = "testjaga6717@gmail.com"
$ws1.Range("D2").Value = "jaga@12345"
$ws1.Range("C3").Value = "testjaga6717@gmail.com"
$ws1.Range("D3").Value = "jaga@1234"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:testjaga6717@gmail.com")
$ws1.Range("C2").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:jaga@12345")
$ws1.Range("D2").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("D3"), "mailto:jaga@1234")
$ws1.Range("D3").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("C3"), "mailto:testjaga6717@gmail.com")
$ws1.Range("C3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet 2 ("CreateAccount"): duplicate row 2 into a new row 3 with a new
# person's data, update row 2's email/phone, then rebuild the hyperlinks.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CreateAccount")

$ws2.Range("A2:O2").Copy($ws2.Range("A3:O3"))

# Update row 2 values per diff (email + phone number change)
$ws2.Range("F2").Value = "testjaga6717@gmail.com"
$ws2.Range("M2").Value = 9878653421

# Fill in row 3's new values
$ws2.Range("C3").Value = "Jack"
$ws2.Range("D3").Value = "Son"
$ws2.Range("E3").Value = "08/30/1993"
$ws2.Range("F3").Value = "testjaga6717@gmail.com"
$ws2.Range("H3").Value = "Mahabalipuram"
$ws2.Range("I3").Value = 600120
$ws2.Range("M3").Value = 9876543211

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("F2"), "mailto:testjaga6717@gmail.com")
$ws2.Range("F2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("N2"), "mailto:jaga@12345")
$ws2.Range("N2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("O2"), "mailto:jaga@12345")
$ws2.Range("O2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("N3"), "mailto:jaga@12345")
$ws2.Range("N3").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("O3"), "mailto:jaga@12345")
$ws2.Range("O3").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("F3"), "mailto:testjaga6717@gmail.com")
$ws2.Range("F3").Style = "Hyperlink"

$ws2.Range("F3").Select()

# ---------------------------------------------------------------------------
# Restore SignIn as the active (tab-selected) sheet with C9 selected, as the
# last action - matches the final UI state captured in the diff.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C9").Select()
